{"js": "// Insert \" on GitHub\" in two places in the \"Present:\" / document-status\n// paragraphs of the assignment:\n//   1. \"...status of your issue\" -> \"...status of your issue on GitHub\"\n//   2. \"...creating this documentation\" -> \"...creating this documentation on GitHub\"\n\nconst body = context.document.body;\n\nconst issueHits = body.search(\"status of your issue\", { matchCase: true });\nissueHits.load(\"items\");\n\nconst docHits = body.search(\"this documentation\", { matchCase: true });\ndocHits.load(\"items\");\n\nawait context.sync();\n\nif (issueHits.items.length > 0) {\n  issueHits.items[0].insertText(\" on GitHub\", \"End\");\n}\n\nif (docHits.items.length > 0) {\n  docHits.items[0].insertText(\" on GitHub\", \"End\");\n}\n\nawait context.sync();\n", "ps1": "# Insert \" on GitHub\" in two places in the \"Present:\" / document-status\n# paragraphs of the assignment:\n#   1. \"...status of your issue\"        -> \"...status of your issue on GitHub\"\n#   2. \"...creating this documentation,\" -> \"...creating this documentation on GitHub,\"\n\n$d = $word.ActiveDocument\n\n$find1 = $d.Content.Find\n$find1.Text = \"status of your issue\"\n$find1.Replacement.Text = \"status of your issue on GitHub\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2) | Out-Null\n\n$find2 = $d.Content.Find\n$find2.Text = \"this documentation,\"\n$find2.Replacement.Text = \"this documentation on GitHub,\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n"}
